$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

# ALC row 70
$wsALC.Range("H70").Value = 6418.364
$wsALC.Range("I70").Value = 26001
$wsALC.Range("J70").Value = 2066.6667
$wsALC.Range("K70").Value = 78003
$wsALC.Range("L70").Value = 6200.000100000001
$wsALC.Range("M70").Value = -77733
$wsALC.Range("N70").Value = -6740.000100000001

# ALC row 73
$wsALC.Range("H73").Value = 6418.364
$wsALC.Range("I73").Value = 26001
$wsALC.Range("J73").Value = 2066.6667
$wsALC.Range("K73").Value = 78003
$wsALC.Range("L73").Value = 6200.000100000001
$wsALC.Range("M73").Value = -77067
$wsALC.Range("N73").Value = -8072.000100000001

# ALC row 100
$wsALC.Range("H100").Value = 2053.6667
$wsALC.Range("I100").Value = 2177.5
$wsALC.Range("J100").Value = 2008.6364
$wsALC.Range("K100").Value = 2177.5
$wsALC.Range("L100").Value = 2008.6364
$wsALC.Range("M100").Value = -1636.5
$wsALC.Range("N100").Value = -3090.6364

# ALC row 137
$wsALC.Range("H137").Value = 3691.8096
$wsALC.Range("I137").Value = 1130.9445
$wsALC.Range("J137").Value = 5612.4585
$wsALC.Range("K137").Value = 3392.8335
$wsALC.Range("L137").Value = 16837.3755
$wsALC.Range("M137").Value = -842.8335000000002
$wsALC.Range("N137").Value = -21937.3755

# ALC row 138
$wsALC.Range("H138").Value = 1897.4906
$wsALC.Range("I138").Value = 1001.45
$wsALC.Range("J138").Value = 2440.5454
$wsALC.Range("K138").Value = 3004.35
$wsALC.Range("L138").Value = 7321.6362
$wsALC.Range("M138").Value = 2135.65
$wsALC.Range("N138").Value = -17601.6362

# ARM row 28
$wsARM.Range("H28").Value = 3920.3333
$wsARM.Range("I28").Value = 2035.375
$wsARM.Range("K28").Value = 2035.375
$wsARM.Range("M28").Value = -1843.375

# ARM row 32
$wsARM.Range("H32").Value = 16403854
$wsARM.Range("I32").Value = 19611476
$wsARM.Range("J32").Value = 44979.9
$wsARM.Range("K32").Value = 19611476
$wsARM.Range("L32").Value = 44979.9
$wsARM.Range("M32").Value = -19611189
$wsARM.Range("N32").Value = -45553.9

# ARM row 81
$wsARM.Range("H81").Value = 40181
$wsARM.Range("J81").Value = 40181
$wsARM.Range("L81").Value = 40181
$wsARM.Range("N81").Value = -42177

# ARM row 84
$wsARM.Range("H84").Value = 40181
$wsARM.Range("J84").Value = 40181
$wsARM.Range("L84").Value = 120543
$wsARM.Range("N84").Value = -130527

# ARM row 95
$wsARM.Range("H95").Value = 20208
$wsARM.Range("J95").Value = 20208
$wsARM.Range("L95").Value = 20208
$wsARM.Range("N95").Value = -25700

# ARM row 97
$wsARM.Range("H97").Value = 883.6111
$wsARM.Range("I97").Value = 709.8182
$wsARM.Range("J97").Value = 1156.7142
$wsARM.Range("K97").Value = 709.8182
$wsARM.Range("L97").Value = 1156.7142
$wsARM.Range("M97").Value = -213.8182
$wsARM.Range("N97").Value = -2148.7142

# ARM row 99
$wsARM.Range("H99").Value = 3920.3333
$wsARM.Range("I99").Value = 2035.375
$wsARM.Range("K99").Value = 2035.375
$wsARM.Range("M99").Value = 959.625

# BSM row 134
$wsBSM.Range("H134").Value = 1838.5238
$wsBSM.Range("I134").Value = 1681.3529
$wsBSM.Range("J134").Value = 2506.5
$wsBSM.Range("K134").Value = 5044.0587
$wsBSM.Range("L134").Value = 7519.5
$wsBSM.Range("M134").Value = -2509.0587
$wsBSM.Range("N134").Value = -12589.5

# CRP row 22
$wsCRP.Range("H22").Value = 190.42857
$wsCRP.Range("I22").Value = 168.66667
$wsCRP.Range("J22").Value = 196.36363
$wsCRP.Range("K22").Value = 168.66667
$wsCRP.Range("L22").Value = 196.36363
$wsCRP.Range("M22").Value = 181.33333
$wsCRP.Range("N22").Value = -896.3636300000001

# CRP row 31
$wsCRP.Range("H31").Value = 3627.32
$wsCRP.Range("I31").Value = 2451.4614
$wsCRP.Range("J31").Value = 4901.1665
$wsCRP.Range("K31").Value = 2451.4614
$wsCRP.Range("L31").Value = 4901.1665
$wsCRP.Range("M31").Value = -2156.4614
$wsCRP.Range("N31").Value = -5491.1665

# CRP row 34
$wsCRP.Range("H34").Value = 3627.32
$wsCRP.Range("I34").Value = 2451.4614
$wsCRP.Range("J34").Value = 4901.1665
$wsCRP.Range("K34").Value = 2451.4614
$wsCRP.Range("L34").Value = 4901.1665
$wsCRP.Range("M34").Value = -2249.4614
$wsCRP.Range("N34").Value = -5305.1665

# CRP row 43
$wsCRP.Range("H43").Value = 0
$wsCRP.Range("J43").Value = 0
$wsCRP.Range("L43").Value = 0
$wsCRP.Range("N43").ClearContents()

# CRP row 92
$wsCRP.Range("H92").Value = 0
$wsCRP.Range("J92").Value = 0
$wsCRP.Range("L92").Value = 0
$wsCRP.Range("N92").ClearContents()

# CRP row 93
$wsCRP.Range("H93").Value = 12916.667
$wsCRP.Range("I93").Value = 10500
$wsCRP.Range("J93").Value = 25000
$wsCRP.Range("K93").Value = 10500
$wsCRP.Range("L93").Value = 25000
$wsCRP.Range("M93").Value = -8628
$wsCRP.Range("N93").Value = -28744

# CRP row 95
$wsCRP.Range("H95").Value = 10737.5
$wsCRP.Range("J95").Value = 10737.5
$wsCRP.Range("L95").Value = 10737.5
$wsCRP.Range("N95").Value = -16229.5

# CRP row 96
$wsCRP.Range("H96").Value = 8946.944
$wsCRP.Range("J96").Value = 8946.944
$wsCRP.Range("L96").Value = 8946.944
$wsCRP.Range("N96").Value = -14438.944

# CRP row 101
$wsCRP.Range("H101").Value = 0
$wsCRP.Range("J101").Value = 0
$wsCRP.Range("L101").Value = 0
$wsCRP.Range("N101").ClearContents()

# CRP row 102
$wsCRP.Range("H102").Value = 0
$wsCRP.Range("J102").Value = 0
$wsCRP.Range("L102").Value = 0
$wsCRP.Range("N102").ClearContents()

# CRP row 103
$wsCRP.Range("H103").Value = 12500
$wsCRP.Range("I103").Value = 12500
$wsCRP.Range("J103").Value = 0
$wsCRP.Range("K103").Value = 12500
$wsCRP.Range("L103").Value = 0
$wsCRP.Range("M103").Value = -11328
$wsCRP.Range("N103").ClearContents()

# CRP row 106
$wsCRP.Range("H106").Value = 19328.2
$wsCRP.Range("J106").Value = 19328.2
$wsCRP.Range("L106").Value = 19328.2
$wsCRP.Range("N106").Value = -21852.2

# CRP row 122
$wsCRP.Range("H122").Value = 1330.7222
$wsCRP.Range("I122").Value = 737.9474
$wsCRP.Range("J122").Value = 1993.2354
$wsCRP.Range("K122").Value = 2213.8422
$wsCRP.Range("L122").Value = 5979.706200000001
$wsCRP.Range("M122").Value = 236.1578
$wsCRP.Range("N122").Value = -10879.7062

# GSM row 132
$wsGSM.Range("H132").Value = 3002.2334
$wsGSM.Range("I132").Value = 2537.8
$wsGSM.Range("J132").Value = 3931.1
$wsGSM.Range("K132").Value = 7613.400000000001
$wsGSM.Range("L132").Value = 11793.3
$wsGSM.Range("M132").Value = -5083.400000000001
$wsGSM.Range("N132").Value = -16853.3

# LTW row 22
$wsLTW.Range("H22").Value = 1196.6666
$wsLTW.Range("I22").Value = 993.3333
$wsLTW.Range("J22").Value = 1400
$wsLTW.Range("K22").Value = 993.3333
$wsLTW.Range("L22").Value = 1400
$wsLTW.Range("M22").Value = -698.3333
$wsLTW.Range("N22").Value = -1990

# LTW row 27
$wsLTW.Range("H27").Value = 1196.6666
$wsLTW.Range("I27").Value = 993.3333
$wsLTW.Range("J27").Value = 1400
$wsLTW.Range("K27").Value = 993.3333
$wsLTW.Range("L27").Value = 1400
$wsLTW.Range("M27").Value = -886.3333
$wsLTW.Range("N27").Value = -1614

# LTW row 40
$wsLTW.Range("H40").Value = 3159.3965
$wsLTW.Range("I40").Value = 2683.375
$wsLTW.Range("J40").Value = 4217.222
$wsLTW.Range("K40").Value = 2683.375
$wsLTW.Range("L40").Value = 4217.222
$wsLTW.Range("M40").Value = -2547.375
$wsLTW.Range("N40").Value = -4489.222

# LTW row 46
$wsLTW.Range("H46").Value = 1637.5
$wsLTW.Range("I46").Value = 1000
$wsLTW.Range("J46").Value = 1728.5714
$wsLTW.Range("K46").Value = 1000
$wsLTW.Range("L46").Value = 1728.5714
$wsLTW.Range("M46").Value = -812
$wsLTW.Range("N46").Value = -2104.5714

# LTW row 88
$wsLTW.Range("H88").Value = 20587.8
$wsLTW.Range("I88").Value = 0
$wsLTW.Range("J88").Value = 20587.8
$wsLTW.Range("K88").Value = 0
$wsLTW.Range("L88").Value = 20587.8
$wsLTW.Range("M88").ClearContents()
$wsLTW.Range("N88").Value = -21443.8

# LTW row 91
$wsLTW.Range("H91").Value = 20587.8
$wsLTW.Range("I91").Value = 0
$wsLTW.Range("J91").Value = 20587.8
$wsLTW.Range("K91").Value = 0
$wsLTW.Range("L91").Value = 20587.8
$wsLTW.Range("M91").ClearContents()
$wsLTW.Range("N91").Value = -23551.8

# LTW row 101
$wsLTW.Range("H101").Value = 5787.3335
$wsLTW.Range("J101").Value = 5787.3335
$wsLTW.Range("L101").Value = 5787.3335
$wsLTW.Range("N101").Value = -12277.3335

# LTW row 102
$wsLTW.Range("H102").Value = 0
$wsLTW.Range("J102").Value = 0
$wsLTW.Range("L102").Value = 0
$wsLTW.Range("N102").ClearContents()

# LTW row 104
$wsLTW.Range("H104").Value = 16611.285
$wsLTW.Range("J104").Value = 16611.285
$wsLTW.Range("L104").Value = 16611.285
$wsLTW.Range("N104").Value = -23599.285

# LTW row 105
$wsLTW.Range("H105").Value = 20000
$wsLTW.Range("J105").Value = 20000
$wsLTW.Range("L105").Value = 20000
$wsLTW.Range("N105").Value = -26988

# WVR row 82
$wsWVR.Range("H82").Value = 0
$wsWVR.Range("J82").Value = 0
$wsWVR.Range("L82").Value = 0
$wsWVR.Range("N82").ClearContents()

# WVR row 85
$wsWVR.Range("H85").Value = 0
$wsWVR.Range("J85").Value = 0
$wsWVR.Range("L85").Value = 0
$wsWVR.Range("N85").ClearContents()

# WVR row 93
$wsWVR.Range("H93").Value = 0
$wsWVR.Range("J93").Value = 0
$wsWVR.Range("L93").Value = 0
$wsWVR.Range("N93").ClearContents()

# WVR row 99
$wsWVR.Range("H99").Value = 0
$wsWVR.Range("J99").Value = 0
$wsWVR.Range("L99").Value = 0
$wsWVR.Range("N99").ClearContents()

# WVR row 122
$wsWVR.Range("H122").Value = 2204.742
$wsWVR.Range("I122").Value = 2257.6
$wsWVR.Range("J122").Value = 2108.6365
$wsWVR.Range("K122").Value = 6772.799999999999
$wsWVR.Range("L122").Value = 6325.9095
$wsWVR.Range("M122").Value = -4322.799999999999
$wsWVR.Range("N122").Value = -11225.9095

# WVR row 126
$wsWVR.Range("H126").Value = 2081.3
$wsWVR.Range("I126").Value = 1652.9166
$wsWVR.Range("J126").Value = 3794.8333
$wsWVR.Range("K126").Value = 4958.7498
$wsWVR.Range("L126").Value = 11384.4999
$wsWVR.Range("M126").Value = -2488.7498
$wsWVR.Range("N126").Value = -16324.4999
